$wb = $excel.ActiveWorkbook

# --- Sheet "Fares": update Child/Senior Peak & Off-Peak prices to mirror Adult prices ---
$fares = $wb.Worksheets.Item("Fares")

# Row 3 (ONE)
$fares.Range("D3").Value = 10
$fares.Range("E3").Value = 8
$fares.Range("F3").Value = 10
$fares.Range("G3").Value = 8

# Row 4 (TWO)
$fares.Range("D4").Value = 20
$fares.Range("E4").Value = 16
$fares.Range("F4").Value = 20
$fares.Range("G4").Value = 16

# Row 5 (THREE)
$fares.Range("D5").Value = 30
$fares.Range("E5").Value = 24
$fares.Range("F5").Value = 30
$fares.Range("G5").Value = 24

# Update the selection shown on the Fares sheet
$fares.Activate()
$fares.Range("G9").Select()

# --- Sheet "Products": update the active selection ---
$products = $wb.Worksheets.Item("Products")
$products.Activate()
$products.Range("F8").Select()

# Restore "Journeys" as the active/visible tab (unchanged by this edit)
$journeys = $wb.Worksheets.Item("Journeys")
$journeys.Activate()
